$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.64"
$ws.Range("E2").Value = "'-0.49%"
$ws.Range("D3").Value = "'37.29"
$ws.Range("E3").Value = "'-0.53%"
$ws.Range("D4").Value = "'5.127"
$ws.Range("E4").Value = "'0.29%"
$ws.Range("D5").Value = "'0.07842"
$ws.Range("E5").Value = "'0.08%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.391"
$ws.Range("E6").Value = "'-0.18%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.257"
$ws.Range("E7").Value = "'0.37%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.882"
$ws.Range("E8").Value = "'-1.83%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.972"
$ws.Range("E9").Value = "'9.04%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9219"
$ws.Range("E10").Value = "'-0.77%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1110"
$ws.Range("E11").Value = "'-7.29%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1903"
$ws.Range("E12").Value = "'0.15%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.08894"
$ws.Range("E13").Value = "'-4.59%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03333"
$ws.Range("E14").Value = "'-2.69%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09604"
$ws.Range("E15").Value = "'-0.09%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001382"
$ws.Range("E16").Value = "'1.26%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005996"
$ws.Range("E17").Value = "'2.69%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.393"
$ws.Range("E18").Value = "'-4.07%"
$ws.Range("D19").Value = "'0.3457"
$ws.Range("E19").Value = "'0.94%"
$ws.Range("D20").Value = "'6.397"
$ws.Range("E20").Value = "'21.59%"
$ws.Range("E21").Value = "'1.80%"
$ws.Range("D22").Value = "'0.2403"
$ws.Range("E22").Value = "'-7.14%"
$ws.Range("D23").Value = "'0.04349"
$ws.Range("E23").Value = "'-0.06%"
$ws.Range("D24").Value = "'0.001200"
$ws.Range("E24").Value = "'0.34%"
$ws.Range("E25").Value = "'0.47%"
$ws.Range("E26").Value = "'7.77%"
$ws.Range("E39").Value = "'4.05%"
$ws.Range("D40").Value = "'0.05025"
$ws.Range("E40").Value = "'-0.42%"
$ws.Range("D41").Value = "'0.007570"
$ws.Range("E41").Value = "'-0.13%"
$ws.Range("D42").Value = "'0.1354"
$ws.Range("E42").Value = "'0.35%"
$ws.Range("D43").Value = "'0.008504"
$ws.Range("E43").Value = "'-6.74%"
$ws.Range("D44").Value = "'0.002070"
$ws.Range("E44").Value = "'3.35%"
$ws.Range("D45").Value = "'0.008137"
$ws.Range("E45").Value = "'-5.47%"
$ws.Range("D46").Value = "'0.00006521"
$ws.Range("E46").Value = "'-2.66%"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.003296"
$ws.Range("E48").Value = "'13.27%"
$ws.Range("E49").Value = "'20.35%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.09%"
